# Apply the "Updated symbol list" price/data refresh to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that are stored as TEXT in the
# workbook (inline strings). Excel auto-converts plain numeric strings to
# real numbers when assigned via .Value, so force a text number format on
# every D cell we touch before writing the new value, keeping them text.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D10","D11","D12","D13","D14",
            "D15","D16","D18","D19","D20","D21","D23","D40","D41","D42","D43",
            "D44","D45","D48","D49")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Simple price (column D) updates -------------------------------------
$ws.Range("D2").Value  = "245.08"
$ws.Range("D3").Value  = "23.06"
$ws.Range("D4").Value  = "5.403"
$ws.Range("D5").Value  = "0.06047"
$ws.Range("D6").Value  = "3.394"
$ws.Range("D7").Value  = "0.8073"
$ws.Range("D8").Value  = "0.9326"
$ws.Range("D10").Value = "0.07431"
$ws.Range("D11").Value = "0.03370"
$ws.Range("D12").Value = "0.03072"
$ws.Range("D13").Value = "0.09363"
$ws.Range("D14").Value = "3.944"
$ws.Range("D15").Value = "0.001587"
$ws.Range("D16").Value = "0.04819"
$ws.Range("D18").Value = "0.005449"
$ws.Range("D19").Value = "0.004166"
$ws.Range("D20").Value = "0.0009864"
$ws.Range("D21").Value = "0.00008706"
$ws.Range("D23").Value = "6.442"
$ws.Range("D40").Value = "0.03977"

# --- Rows 41-43: coins reshuffled (Kick -> BKEX -> CEJI -> Kick) ----------
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1076"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002712"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.003037"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- Remaining price / label updates --------------------------------------
$ws.Range("D44").Value = "0.005947"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"

$ws.Range("D45").Value = "0.00005201"

$ws.Range("D48").Value = "0.8204"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"

$ws.Range("D49").Value = "0.002179"
